$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9446313977241516
$ws.Range("B1").Value = 1.682584881782532
$ws.Range("C1").Value = 5.557644844055176
$ws.Range("D1").Value = 3.593063831329346
$ws.Range("E1").Value = 1.4019855260849
